# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1 (AD, AE, AF) - same bold/border/center style as
# the rest of the header row (A1:AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats

# Every player row (2-56) gets the same season team record.
$lastRow = 56
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 98   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 64   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
